# Commit: "Fruta / hortaliza, semanal"
# A new weekly data row is inserted at row 383 (pushing the existing rows
# 383-471 down to 384-472). The workbook only has one worksheet, which is
# the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 383, shifting all rows below it
# down by one (this also grows the sheet dimension from R471 to R472).
$ws.Rows("383:383").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A383").Value = 4
$ws.Range("B383").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C383").Value = "Los Lagos"
$ws.Range("D383").Value = 45275
$ws.Range("E383").Value = 10
$ws.Range("F383").Value = 100112028
$ws.Range("G383").Value = "Sandia"
$ws.Range("H383").Value = "Sin especificar"
$ws.Range("I383").Value = "Primera"
$ws.Range("J383").Value = 800
$ws.Range("K383").Value = 1000
$ws.Range("L383").Value = 1200
$ws.Range("M383").Value = 1100
$ws.Range("N383").Value = "$/kilo (volumen en unidades)"
$ws.Range("O383").Value = "Perú"
$ws.Range("P383").Value = 1100
$ws.Range("Q383").Value = 1
$ws.Range("R383").Value = "Hortaliza"
